# Applies the edits described by the CU_Enviar_Mensaje_Privado.docx diff:
#  1. "Caso de Uso: Visualizar Perfiles de Usuarios" -> "Caso de Uso: Ver Perfil Usuario"
#     (split into separate runs: "Ver Perfil" / " " / "Usuario")
#  2. "Hacer click en el boton "Mensajes privados" en la parte superior del Microblog."
#     gets re-run-split with proofing-error markers around "click" (grammar) and
#     "Microblog" (spelling), matching Word's automatic proofing markup.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Edit 1: "Visualizar Perfiles de Usuarios" -> "Ver Perfil Usuario"
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$found1 = $rng1.Find.Execute("Caso de Uso: Visualizar Perfiles de Usuarios", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Edit 1: target paragraph text not found"
}

$xml1 = @'
<?xml version="1.0" encoding="UTF-8"?><w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Caso de Uso: </w:t></w:r><w:r><w:t>Ver Perfil</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>Usuario</w:t></w:r></w:p></w:body></w:wordDocument>
'@
$rng1.InsertXML($xml1)
Write-Output "Edit 1 applied"

# ---------------------------------------------------------------------------
# Edit 2: "Hacer click en el boton..." paragraph gains proofErr markup
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$target2 = "Hacer click en el bot" + [char]0x00F3 + "n en el " + [char]0x201C + "Mensajes privados" + [char]0x201D
# Use the full unique sentence (quotes/accents rebuilt from code points so the
# script source encoding can't bite us) to find the exact paragraph.
$needle2 = "Hacer click en el bot" + [char]0x00F3 + "n " + [char]0x201C + "Mensajes privados" + [char]0x201D + " en la parte superior del Microblog."
$found2 = $rng2.Find.Execute($needle2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Edit 2: target paragraph text not found"
}

$xml2 = @'
<?xml version="1.0" encoding="UTF-8"?><w:wordDocument xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:widowControl w:val="0"/><w:spacing w:line="240" w:lineRule="auto"/></w:pPr><w:r><w:t xml:space="preserve">Hacer </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>click</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> en el bot&#243;n &#8220;Mensajes privados&#8221; en la parte superior del </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Microblog</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>.</w:t></w:r></w:p></w:body></w:wordDocument>
'@
$rng2.InsertXML($xml2)
Write-Output "Edit 2 applied"

Write-Output "Done"
